$d = $word.ActiveDocument

# 1) Capitalize "django" -> "Django" in "the django debug on"
$d.Content.Find.Execute("the django debug on", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "the Django debug on", 2)

# 2) Replace the trailing clause with the new, longer explanation
$d.Content.Find.Execute(", which allows the users to see our URL patterns.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, " when debug is off, leading to Django debug pages if the user gets the site to crash. We added all sorts of checks to prevent crashes from being possible.", 2)
